{"js": "// Update the worksheet date and all twenty-five \"three-digit \u00f7 one-digit\"\n// division prompts to the next day's generated set of problems.\nconst replacements = [\n  [\"2024-11-11 Monday\", \"2024-11-12 Tuesday\"],\n  [\"401\u00f73=\", \"814\u00f79=\"],\n  [\"566\u00f72=\", \"430\u00f73=\"],\n  [\"855\u00f74=\", \"462\u00f74=\"],\n  [\"126\u00f72=\", \"945\u00f74=\"],\n  [\"600\u00f75=\", \"742\u00f75=\"],\n  [\"943\u00f76=\", \"339\u00f73=\"],\n  [\"498\u00f78=\", \"438\u00f75=\"],\n  [\"157\u00f72=\", \"219\u00f76=\"],\n  [\"926\u00f76=\", \"704\u00f73=\"],\n  [\"821\u00f72=\", \"287\u00f73=\"],\n  [\"964\u00f78=\", \"868\u00f73=\"],\n  [\"175\u00f72=\", \"824\u00f75=\"],\n  [\"439\u00f79=\", \"544\u00f76=\"],\n  [\"680\u00f76=\", \"130\u00f73=\"],\n  [\"260\u00f75=\", \"606\u00f72=\"],\n  [\"448\u00f79=\", \"733\u00f76=\"],\n  [\"457\u00f77=\", \"803\u00f78=\"],\n  [\"796\u00f79=\", \"217\u00f79=\"],\n  [\"198\u00f77=\", \"141\u00f79=\"],\n  [\"125\u00f75=\", \"837\u00f77=\"],\n  [\"604\u00f77=\", \"924\u00f77=\"],\n  [\"977\u00f77=\", \"425\u00f78=\"],\n  [\"921\u00f75=\", \"810\u00f74=\"],\n  [\"769\u00f76=\", \"391\u00f76=\"],\n  [\"188\u00f74=\", \"982\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all twenty-five \"three-digit \u00f7 one-digit\"\n# division prompts to the next day's generated set of problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-11 Monday\", \"2024-11-12 Tuesday\"),\n    @(\"401\u00f73=\", \"814\u00f79=\"),\n    @(\"566\u00f72=\", \"430\u00f73=\"),\n    @(\"855\u00f74=\", \"462\u00f74=\"),\n    @(\"126\u00f72=\", \"945\u00f74=\"),\n    @(\"600\u00f75=\", \"742\u00f75=\"),\n    @(\"943\u00f76=\", \"339\u00f73=\"),\n    @(\"498\u00f78=\", \"438\u00f75=\"),\n    @(\"157\u00f72=\", \"219\u00f76=\"),\n    @(\"926\u00f76=\", \"704\u00f73=\"),\n    @(\"821\u00f72=\", \"287\u00f73=\"),\n    @(\"964\u00f78=\", \"868\u00f73=\"),\n    @(\"175\u00f72=\", \"824\u00f75=\"),\n    @(\"439\u00f79=\", \"544\u00f76=\"),\n    @(\"680\u00f76=\", \"130\u00f73=\"),\n    @(\"260\u00f75=\", \"606\u00f72=\"),\n    @(\"448\u00f79=\", \"733\u00f76=\"),\n    @(\"457\u00f77=\", \"803\u00f78=\"),\n    @(\"796\u00f79=\", \"217\u00f79=\"),\n    @(\"198\u00f77=\", \"141\u00f79=\"),\n    @(\"125\u00f75=\", \"837\u00f77=\"),\n    @(\"604\u00f77=\", \"924\u00f77=\"),\n    @(\"977\u00f77=\", \"425\u00f78=\"),\n    @(\"921\u00f75=\", \"810\u00f74=\"),\n    @(\"769\u00f76=\", \"391\u00f76=\"),\n    @(\"188\u00f74=\", \"982\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
